# Repo cleanup before site publish: outlet menu link updates + engine implementation.
# This script converts the "salads_greenmountain" sheet into "pokebowls_greenmountain":
#  - renames the sheet
#  - fills in the missing "Allergens" value for the Tuna row
#  - extends the table with a new (currently empty/header-only) "Column1" column
#  - updates the table style
#  - normalizes row heights / removes the now-unused wrap-text style

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the new menu.
$ws.Name = "pokebowls_greenmountain"

# Fill in the previously-blank Allergens cell for the Tuna row.
$ws.Range("C4").Value = "No known priority allergens"

# Remove the wrap-text formatting that was only used by the
# Ingredients/Allergens cells of the Coastal & Veggie rows.
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Style = "Normal"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Style = "Normal"

# Normalize row heights to a compact 12pt custom height across the table.
$ws.Rows.Item(1).RowHeight = 12
$ws.Rows.Item(2).RowHeight = 12
$ws.Rows.Item(3).RowHeight = 12
$ws.Rows.Item(4).RowHeight = 12
$ws.Rows.Item(5).RowHeight = 12
$ws.Rows.Item(6).RowHeight = 12

# Extend the table to include a new trailing column and give it a header.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G6"))
$ws.Range("G1").Value = "Column1"

# Apply the new table style.
$lo.TableStyle = "TableStyleMedium7"

# Update the active selection to match the published state.
$ws.Range("C4").Select()
